# Slide 1 ("Title" / intro slide): reposition three shapes, matching the
# author's drag-and-drop move of the title placeholder, the bullet-list
# "Title 1" textbox, and the "2" rectangle callout.
#
# Shape.Left/Top are expressed in points (1 pt = 12700 EMU) and are stored
# internally as single-precision floats, so the literals below are nudged to
# the nearest value whose float32 representation still truncates to the
# exact target EMU offset from the target OOXML (EMU = floor(Single(pt) * 12700)).
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# "Title 1" (ctrTitle placeholder, id=12290): (2819400,1981200) -> (2673096,1585551) EMU
$s.Shapes.Item(1).Left = 210.48001098632812
$s.Shapes.Item(1).Top = 124.84654235839844

# "Title 1" (bullet textbox, id=12292): (3962400,4352330) -> (2523744,4335334) EMU
$s.Shapes.Item(3).Left = 198.72000122070312
$s.Shapes.Item(3).Top = 341.3648986816406

# "Rectangle 4" (id=5): (8010144,2658934) -> (7918704,2246288) EMU
$s.Shapes.Item(4).Left = 623.52001953125
$s.Shapes.Item(4).Top = 176.87307739257812
